# C5-PowerPoint.pptx - Fri, Jul 03, 2020  9:05:33 AM
#
# 1) The table on slide 6 is re-styled with a different (built-in) table
#    style id.
# 2) The colour palette that backs the deck's design theme is swapped from
#    the "Integral" palette to the "Office Theme" palette (the theme's font
#    scheme / format scheme were already identical between the two themes,
#    so only the twelve theme colours actually change visually).

$p = $ppt.ActivePresentation

# --- 1. Table re-style (slide 6) -------------------------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shape = $tableSlide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{1700A436-AC7D-4B8A-B49F-957A0A8F68A3}")
    }
}

# --- 2. Theme colour swap (Integral -> Office Theme) -----------------------
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    0,         # dk1      #000000
    16777215,  # lt1      #FFFFFF
    6968388,   # dk2      #44546A
    15132391,  # lt2      #E7E6E6
    13998939,  # accent1  #5B9BD5
    3243501,   # accent2  #ED7D31
    10855845,  # accent3  #A5A5A5
    49407,     # accent4  #FFC000
    12874308,  # accent5  #4472C4
    4697456,   # accent6  #70AD47
    12673797,  # hlink    #0563C1
    7491477    # folHlink #954F72
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
